# Auto-generated Excel COM-interop script applying the FlashScore odds
# refresh described by the "Atualizando o arquivo XLSX" commit.
#
# The workbook has a single worksheet (Sheet1) listing football matches
# for 2025-05-14 together with betting odds / count columns (G:AJ).
# The edit updates a set of individual odds/count cells to their new
# values (re-scraped odds), and fills in odds data for two matches
# (rows 36-37) that previously had empty odds columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Almere City - Sittard
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 3.2
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 4.33
$ws.Range("N3").Value = 1.7
$ws.Range("O3").Value = 2.1
$ws.Range("R3").Value = 1.57
$ws.Range("S3").Value = 2.25
$ws.Range("AA3").Value = 7
$ws.Range("AG3").Value = 11

# Row 4: Feyenoord - Waalwijk
$ws.Range("N4").Value = 1.29
$ws.Range("O4").Value = 3.6

# Row 6: Groningen - Ajax
$ws.Range("L6").Value = 1.2
$ws.Range("M6").Value = 4.33
$ws.Range("N6").Value = 1.67
$ws.Range("O6").Value = 2.15

# Row 8: PSV - Heracles
$ws.Range("G8").Value = 1.07
$ws.Range("H8").Value = 12
$ws.Range("R8").Value = 2.05
$ws.Range("S8").Value = 1.7
$ws.Range("W8").Value = 7.5

# Row 11: Willem II - Zwolle
$ws.Range("G11").Value = 2.63
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 13
$ws.Range("N11").Value = 1.83
$ws.Range("O11").Value = 2.03

# Row 13: Villarreal - Leganes
$ws.Range("T13").Value = 8
$ws.Range("Y13").Value = 23
$ws.Range("AE13").Value = 23
$ws.Range("AG13").Value = 23

# Row 14: Real Madrid - Mallorca
$ws.Range("G14").Value = 1.38

# Row 15: Lok. Sofia - CSKA 1948 Sofia
$ws.Range("H15").Value = 3.3
$ws.Range("I15").Value = 3.2
$ws.Range("R15").Value = 1.83
$ws.Range("S15").Value = 1.83
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 15
$ws.Range("AG15").Value = 12

# Row 16: Arda - Levski Sofia
$ws.Range("I16").Value = 2.4
$ws.Range("J16").Value = 1.1
$ws.Range("K16").Value = 7
$ws.Range("AH16").Value = 23

# Row 17: Beroe - Botev Plovdiv
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 2.38
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10
$ws.Range("T17").Value = 10
$ws.Range("U17").Value = 15
$ws.Range("V17").Value = 11
$ws.Range("W17").Value = 29
$ws.Range("Y17").Value = 29
$ws.Range("AE17").Value = 9
$ws.Range("AH17").Value = 23

# Row 18: Alajuelense - Puntarenas FC
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11
$ws.Range("N18").Value = 1.93
$ws.Range("O18").Value = 1.88

# Row 22: Stockport County - Leyton Orient
$ws.Range("G22").Value = 1.8
$ws.Range("H22").Value = 3.7
$ws.Range("I22").Value = 4.1
$ws.Range("P22").Value = 1.36
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 1.8
$ws.Range("S22").Value = 1.95
$ws.Range("T22").Value = 7.5
$ws.Range("Y22").Value = 26
$ws.Range("Z22").Value = 11
$ws.Range("AD22").Value = 251
$ws.Range("AE22").Value = 12

# Row 24: Levadiakos - Panetolikos
$ws.Range("G24").Value = 1.83
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 4.75
$ws.Range("U24").Value = 7.5
$ws.Range("W24").Value = 15
$ws.Range("X24").Value = 19
$ws.Range("Y24").Value = 41

# Row 26: Volos - Lamia
$ws.Range("G26").Value = 1.17
$ws.Range("H26").Value = 7.5
$ws.Range("I26").Value = 19
$ws.Range("V26").Value = 11
$ws.Range("W26").Value = 6
$ws.Range("Z26").Value = 12
$ws.Range("AA26").Value = 15
$ws.Range("AB26").Value = 41
$ws.Range("AC26").Value = 151
$ws.Range("AE26").Value = 29
$ws.Range("AH26").Value = 351
$ws.Range("AI26").Value = 151
$ws.Range("AJ26").Value = 126

# Row 27: FeralpiSalo - Crotone
$ws.Range("K27").Value = 7.1
$ws.Range("L27").Value = 1.3
$ws.Range("M27").Value = 3.25
$ws.Range("N27").Value = 1.87
$ws.Range("O27").Value = 1.83
$ws.Range("P27").Value = 1.44
$ws.Range("Q27").Value = 2.62
$ws.Range("R27").Value = 1.7
$ws.Range("S27").Value = 2.05
$ws.Range("T27").Value = 8
$ws.Range("X27").Value = 16.5
$ws.Range("Z27").Value = 7.1
$ws.Range("AA27").Value = 6.1
$ws.Range("AE27").Value = 10

# Row 30: Rimini - Vis Pesaro
$ws.Range("H30").Value = 3
$ws.Range("I30").Value = 2.95
$ws.Range("Q30").Value = 2.32
$ws.Range("T30").Value = 6.5
$ws.Range("U30").Value = 10.75
$ws.Range("V30").Value = 9.75
$ws.Range("X30").Value = 24
$ws.Range("Y30").Value = 40
$ws.Range("AA30").Value = 5.8
$ws.Range("AE30").Value = 7.8
$ws.Range("AF30").Value = 14.5
$ws.Range("AG30").Value = 10.75
$ws.Range("AI30").Value = 28
$ws.Range("AJ30").Value = 40

# Row 31: Torres - Atalanta U23
$ws.Range("G31").Value = 2.32
$ws.Range("H31").Value = 3.15
$ws.Range("I31").Value = 3
$ws.Range("J31").Value = 1.06
$ws.Range("U31").Value = 11.75
$ws.Range("V31").Value = 9
$ws.Range("W31").Value = 24
$ws.Range("X31").Value = 18.5
$ws.Range("Y31").Value = 27
$ws.Range("AB31").Value = 12.5
$ws.Range("AF31").Value = 16
$ws.Range("AG31").Value = 10.5
$ws.Range("AH31").Value = 37
$ws.Range("AI31").Value = 25
$ws.Range("AJ31").Value = 30

# Row 33: Yokohama F. Marinos - Kashiwa Reysol
$ws.Range("N33").Value = 2.03
$ws.Range("O33").Value = 1.83

# Row 36: Dep. Capiata - Guarani de Fram
$ws.Range("G36").Value = 1.75
$ws.Range("H36").Value = 3.45
$ws.Range("I36").Value = 4.4
$ws.Range("L36").Value = 1.33
$ws.Range("M36").Value = 2.8
$ws.Range("N36").Value = 1.98
$ws.Range("O36").Value = 1.65
$ws.Range("P36").Value = 1.44
$ws.Range("Q36").Value = 2.42
$ws.Range("R36").Value = 1.88
$ws.Range("S36").Value = 1.72
$ws.Range("T36").Value = 6.1
$ws.Range("U36").Value = 7.6
$ws.Range("V36").Value = 8.25
$ws.Range("W36").Value = 13.5
$ws.Range("X36").Value = 15
$ws.Range("Y36").Value = 32
$ws.Range("Z36").Value = 8.75
$ws.Range("AA36").Value = 6.7
$ws.Range("AB36").Value = 17
$ws.Range("AC36").Value = 90
$ws.Range("AD36").Value = 800
$ws.Range("AE36").Value = 11.25
$ws.Range("AF36").Value = 24
$ws.Range("AG36").Value = 14.5
$ws.Range("AH36").Value = 75
$ws.Range("AI36").Value = 45
$ws.Range("AJ36").Value = 55

# Row 37: Independiente FBC - Tacuary
$ws.Range("G37").Value = 2.52
$ws.Range("H37").Value = 3.2
$ws.Range("I37").Value = 2.6
$ws.Range("L37").Value = 1.38
$ws.Range("M37").Value = 2.57
$ws.Range("N37").Value = 2.12
$ws.Range("O37").Value = 1.57
$ws.Range("P37").Value = 1.47
$ws.Range("Q37").Value = 2.35
$ws.Range("R37").Value = 1.88
$ws.Range("S37").Value = 1.72
$ws.Range("T37").Value = 7.1
$ws.Range("U37").Value = 11.5
$ws.Range("V37").Value = 10
$ws.Range("W37").Value = 27
$ws.Range("X37").Value = 23
$ws.Range("Y37").Value = 40
$ws.Range("Z37").Value = 8
$ws.Range("AA37").Value = 6.2
$ws.Range("AB37").Value = 16.5
$ws.Range("AC37").Value = 90
$ws.Range("AD37").Value = 900
$ws.Range("AE37").Value = 7.3
$ws.Range("AF37").Value = 11.75
$ws.Range("AG37").Value = 10.25
$ws.Range("AH37").Value = 28
$ws.Range("AI37").Value = 24
$ws.Range("AJ37").Value = 40

# Row 39: Leixoes - Feirense
$ws.Range("L39").Value = 1.29
$ws.Range("M39").Value = 3.5
$ws.Range("N39").Value = 1.9
$ws.Range("O39").Value = 1.9

# Row 41: Hearts - St Johnstone
$ws.Range("G41").Value = 1.65
$ws.Range("H41").Value = 3.75
$ws.Range("I41").Value = 5.25
$ws.Range("AA41").Value = 7
$ws.Range("AB41").Value = 15
$ws.Range("AG41").Value = 17

# Row 44: St. Mirren - Hibernian
$ws.Range("G44").Value = 2.5
$ws.Range("I44").Value = 2.75
$ws.Range("N44").Value = 1.98
$ws.Range("O44").Value = 1.88
$ws.Range("T44").Value = 8.5
$ws.Range("U44").Value = 12
$ws.Range("V44").Value = 10
$ws.Range("W44").Value = 23
$ws.Range("X44").Value = 21
$ws.Range("AE44").Value = 9.5
$ws.Range("AF44").Value = 13
$ws.Range("AG44").Value = 10
$ws.Range("AH44").Value = 29
$ws.Range("AI44").Value = 21

# Row 45: Aberdeen - Celtic
$ws.Range("G45").Value = 5.25
$ws.Range("I45").Value = 1.62
$ws.Range("R45").Value = 1.62
$ws.Range("S45").Value = 2.2
$ws.Range("Y45").Value = 34
$ws.Range("AF45").Value = 9
$ws.Range("AH45").Value = 13

# Row 47: Lion City - Balestier Khalsa
$ws.Range("G47").Value = 1.7
$ws.Range("I47").Value = 3.25
$ws.Range("R47").Value = 1.19
$ws.Range("S47").Value = 3.75
$ws.Range("W47").Value = 23
$ws.Range("AG47").Value = 17
$ws.Range("AH47").Value = 41
$ws.Range("AI47").Value = 23
$ws.Range("AJ47").Value = 19

# Row 53: CF Montreal - Columbus Crew
$ws.Range("G53").Value = 3.2
$ws.Range("H53").Value = 3.4
$ws.Range("I53").Value = 2.2
$ws.Range("K53").Value = 12
$ws.Range("P53").Value = 1.33
$ws.Range("Q53").Value = 3.25
$ws.Range("R53").Value = 1.67
$ws.Range("S53").Value = 2.1
$ws.Range("T53").Value = 11
$ws.Range("U53").Value = 17
$ws.Range("AD53").Value = 151
$ws.Range("AJ53").Value = 23

# Row 54: DC United - New York City
$ws.Range("N54").Value = 1.65
$ws.Range("O54").Value = 2.2

# Row 65: San Jose Earthquakes - Inter Miami
$ws.Range("G65").Value = 2.5
$ws.Range("I65").Value = 2.5
$ws.Range("T65").Value = 15
